$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gerber_URL")
$ws.Name = "URL"
$ws.Range("F32").Select()
